$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(44376, 1, 5, 12.41372461393316),
  @(44377, 0, 4, 9.930979691146533),
  @(44378, 2, 5, 12.41372461393316),
  @(44379, 2, 7, 17.37921445950643),
  @(44380, 1, 7, 17.37921445950643),
  @(44381, 2, 8, 19.86195938229307),
  @(44382, 2, 10, 24.82744922786633),
  @(44383, 0, 9, 22.34470430507969),
  @(44384, 0, 9, 22.34470430507969),
  @(44385, 3, 10, 24.82744922786633),
  @(44386, 1, 9, 22.34470430507969),
  @(44387, 0, 8, 19.86195938229307),
  @(44388, 7, 13, 32.27568399622623),
  @(44389, 1, 12, 29.79293907343959),
  @(44390, 2, 14, 34.75842891901286),
  @(44391, 3, 17, 42.20666368737276),
  @(44392, 3, 17, 42.20666368737276),
  @(44393, 2, 18, 44.68940861015939),
  @(44394, 4, 22, 54.62038830130592),
  @(44395, 3, 18, 44.68940861015939),
  @(44396, 8, 25, 62.06862306966583),
  @(44397, 0, 23, 57.10313322409256),
  @(44398, 4, 24, 59.58587814687919),
  @(44399, 3, 24, 59.58587814687919),
  @(44400, 2, 24, 59.58587814687919),
  @(44401, 2, 22, 54.62038830130592),
  @(44402, 2, 21, 52.13764337851929)
)

$startRow = 302
$r = $startRow
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}

$lastRow = $r - 1
$ws.Range("A301").Copy()
$ws.Range("A" + $startRow + ":A" + $lastRow).PasteSpecial(-4122)
